$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 21, shifting existing rows 21-98 down to 22-99.
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with its data (matches the surrounding
# "Repollo" / "Crespo record" / "Primera" records for this market).
$ws.Cells.Item(21, 1).Value = 7
$ws.Cells.Item(21, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(21, 3).Value = "Ñuble"
$ws.Cells.Item(21, 4).Value = 44414
$ws.Cells.Item(21, 5).Value = 16
$ws.Cells.Item(21, 6).Value = 100112006
$ws.Cells.Item(21, 7).Value = "Repollo"
$ws.Cells.Item(21, 8).Value = "Crespo record"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 300
$ws.Cells.Item(21, 11).Value = 600
$ws.Cells.Item(21, 12).Value = 650
$ws.Cells.Item(21, 13).Value = 625
$ws.Cells.Item(21, 14).Value = "$/unidad"
$ws.Cells.Item(21, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(21, 16).Value = 625
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = "Hortaliza"
